$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 values are removed (cleared), C2 and E2 get new values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -4.5192477786255836
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -1.8852598986136755

# Row 3: update values in B3:E3
$ws.Range("B3").Value = -5.4145561567021687
$ws.Range("C3").Value = 3.2813695555772853
$ws.Range("D3").Value = -2.5076467871384907
$ws.Range("E3").Value = 10.568140331805843

# Update selection to reflect the new selected range B1:E3
$ws.Range("B1:E3").Select()
